$d = $word.ActiveDocument

# Reverted magenta text back to yellow, and added the clause noting its
# symbolism (full moon), as described in the commit message.
$d.Content.Find.Execute(
    "magenta colored text to point",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "yellow colored text (which is symbolic for a full moon) to point",
    2
)
